# Swap the weekly values between row 2 and row 4 for columns
# D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion),
# S (Precio $/Kg) and T (Kg / unidad).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "N", "O", "P", "Q", "S", "T")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell4 = $ws.Range($col + "4")

    $val2 = $cell2.Value2
    $val4 = $cell4.Value2

    $cell2.Value2 = $val4
    $cell4.Value2 = $val2
}
